# bom fixed (510k resistor -> 510 resistor)
#
# The sheet had two separate resistor rows that should have been one:
#   row 20: Qty=1  Value=510k  Parts="R49"
#   row 23: Qty=3  Value=510   Parts="R9, R10, R37"
# The "510k" row was really a mis-entered duplicate of R49, which is
# actually a 510-ohm resistor like R9/R10/R37. Fix: remove the bogus
# "510k" row and fold R49 into the 510-ohm row's part list (bumping its
# quantity from 3 to 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole "510k" / R49 row — everything below shifts up by one.
$ws.Rows.Item(20).Delete() | Out-Null

# The 510-ohm resistor row (previously row 23) is now row 22.
# Merge R49 into its Parts list and bump Qty from 3 to 4.
$ws.Range("A22").Value = 4
$ws.Range("E22").Value = "R9, R10, R37, R49"

# Match the author's final selection position.
$ws.Range("A23").Select() | Out-Null
